# Scheduled-runner update: refresh cached Universalis price/profit figures
# across the per-job "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each row holds cached numeric snapshots (no live formulas) for a leve --
# currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) -- so
# the refresh just overwrites the affected cells with the newly fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1339.1333
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 1330.6666
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 3991.9998
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -4327.9998
$ws.Range("H69").Value = 16666.334
$ws.Range("J69").Value = 16666.334
$ws.Range("L69").Value = 49999.00199999999
$ws.Range("N69").Value = -51747.00199999999
$ws.Range("H72").Value = 16666.334
$ws.Range("J72").Value = 16666.334
$ws.Range("L72").Value = 149997.006
$ws.Range("N72").Value = -158733.006
$ws.Range("H100").Value = 2043.875
$ws.Range("I100").Value = 1995.8572
$ws.Range("J100").Value = 2081.2222
$ws.Range("K100").Value = 1995.8572
$ws.Range("L100").Value = 2081.2222
$ws.Range("M100").Value = -1454.8572
$ws.Range("N100").Value = -3163.2222
$ws.Range("H111").Value = 1636.3
$ws.Range("I111").Value = 985.63635
$ws.Range("K111").Value = 2956.90905
$ws.Range("M111").Value = 110.0909499999998
$ws.Range("H132").Value = 4491.157
$ws.Range("I132").Value = 2341.2373
$ws.Range("K132").Value = 7023.711899999999
$ws.Range("M132").Value = -4493.711899999999
$ws.Range("H138").Value = 2145.9507
$ws.Range("I138").Value = 1677.0883
$ws.Range("K138").Value = 5031.2649
$ws.Range("M138").Value = 108.7350999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3953.6667
$ws.Range("J74").Value = 3953.6667
$ws.Range("L74").Value = 3953.6667
$ws.Range("N74").Value = -5701.6667
$ws.Range("H77").Value = 3953.6667
$ws.Range("J77").Value = 3953.6667
$ws.Range("L77").Value = 19768.3335
$ws.Range("N77").Value = -28504.3335
$ws.Range("H139").Value = 107994.11
$ws.Range("J139").Value = 107994.11
$ws.Range("L139").Value = 107994.11
$ws.Range("N139").Value = -118274.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1982.4615
$ws.Range("I105").Value = 2002.5
$ws.Range("K105").Value = 2002.5
$ws.Range("M105").Value = -255.5
$ws.Range("H132").Value = 122400.57
$ws.Range("J132").Value = 122400.57
$ws.Range("L132").Value = 122400.57
$ws.Range("N132").Value = -132520.57
$ws.Range("H137").Value = 99578.5
$ws.Range("J137").Value = 99578.5
$ws.Range("L137").Value = 99578.5
$ws.Range("N137").Value = -109778.5
$ws.Range("H141").Value = 83982.55
$ws.Range("I141").Value = 58990
$ws.Range("K141").Value = 58990
$ws.Range("M141").Value = -53810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1040940.4
$ws.Range("I6").Value = 1300675.5
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 1300675.5
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -1300562.5
$ws.Range("N6").Value = -2226
$ws.Range("H31").Value = 3277.7046
$ws.Range("I31").Value = 1873.1818
$ws.Range("J31").Value = 4682.227
$ws.Range("K31").Value = 1873.1818
$ws.Range("L31").Value = 4682.227
$ws.Range("M31").Value = -1578.1818
$ws.Range("N31").Value = -5272.227
$ws.Range("H34").Value = 3277.7046
$ws.Range("I34").Value = 1873.1818
$ws.Range("J34").Value = 4682.227
$ws.Range("K34").Value = 1873.1818
$ws.Range("L34").Value = 4682.227
$ws.Range("M34").Value = -1671.1818
$ws.Range("N34").Value = -5086.227
$ws.Range("H99").Value = 10052.366
$ws.Range("I99").Value = 6905.909
$ws.Range("K99").Value = 6905.909
$ws.Range("M99").Value = -5407.909
$ws.Range("H126").Value = 10052.366
$ws.Range("I126").Value = 6905.909
$ws.Range("K126").Value = 20717.727
$ws.Range("M126").Value = -18247.727
$ws.Range("H132").Value = 7171.884
$ws.Range("I132").Value = 2459.2
$ws.Range("J132").Value = 13717.277
$ws.Range("K132").Value = 7377.599999999999
$ws.Range("L132").Value = 41151.831
$ws.Range("M132").Value = -4847.599999999999
$ws.Range("N132").Value = -46211.831
$ws.Range("H133").Value = 82993.17999999999
$ws.Range("J133").Value = 85042.5
$ws.Range("L133").Value = 85042.5
$ws.Range("N133").Value = -90102.5
$ws.Range("H141").Value = 227857
$ws.Range("J141").Value = 233846
$ws.Range("L141").Value = 233846
$ws.Range("N141").Value = -244206

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 675.5
$ws.Range("I44").Value = 220.6
$ws.Range("K44").Value = 661.8
$ws.Range("M44").Value = -263.8
$ws.Range("H68").Value = 1462.5714
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1462.5714
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4387.7142
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6009.7142
$ws.Range("H71").Value = 1462.5714
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1462.5714
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 13163.1426
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -21275.1426
$ws.Range("H99").Value = 3186.4285
$ws.Range("I99").Value = 1169.8
$ws.Range("K99").Value = 3509.4
$ws.Range("M99").Value = -1263.4
$ws.Range("H131").Value = 1226823.9
$ws.Range("I131").Value = 1730757.4
$ws.Range("K131").Value = 5192272.199999999
$ws.Range("M131").Value = -5187232.199999999
$ws.Range("H134").Value = 2798.842
$ws.Range("I134").Value = 1655.5714
$ws.Range("K134").Value = 4966.7142
$ws.Range("M134").Value = 103.2857999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 79999
$ws.Range("I24").Value = 79999
$ws.Range("K24").Value = 79999
$ws.Range("M24").Value = -79826
$ws.Range("H80").Value = 36926330
$ws.Range("I80").Value = 60003030
$ws.Range("J80").Value = 3607.6
$ws.Range("K80").Value = 60003030
$ws.Range("L80").Value = 3607.6
$ws.Range("M80").Value = -60002032
$ws.Range("N80").Value = -5603.6
$ws.Range("H83").Value = 36926330
$ws.Range("I83").Value = 60003030
$ws.Range("J83").Value = 3607.6
$ws.Range("K83").Value = 300015150
$ws.Range("L83").Value = 18038
$ws.Range("M83").Value = -300010158
$ws.Range("N83").Value = -28022
$ws.Range("H86").Value = 54020
$ws.Range("J86").Value = 54020
$ws.Range("L86").Value = 54020
$ws.Range("N86").Value = -56392
$ws.Range("H89").Value = 54020
$ws.Range("J89").Value = 54020
$ws.Range("L89").Value = 162060
$ws.Range("N89").Value = -173916
$ws.Range("H132").Value = 2422.7144
$ws.Range("I132").Value = 1600.3636
$ws.Range("K132").Value = 4801.0908
$ws.Range("M132").Value = -2271.0908
$ws.Range("H135").Value = 134158.33
$ws.Range("J135").Value = 134158.33
$ws.Range("L135").Value = 134158.33
$ws.Range("N135").Value = -144298.33
$ws.Range("H141").Value = 44723.816
$ws.Range("J141").Value = 44723.816
$ws.Range("L141").Value = 44723.816
$ws.Range("N141").Value = -55083.816

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5828.4287
$ws.Range("I7").Value = 3701.3333
$ws.Range("K7").Value = 3701.3333
$ws.Range("M7").Value = -3589.3333
$ws.Range("H9").Value = 227.25
$ws.Range("I9").Value = 273.33334
$ws.Range("J9").Value = 89
$ws.Range("K9").Value = 273.33334
$ws.Range("L9").Value = 89
$ws.Range("M9").Value = -49.33334000000002
$ws.Range("N9").Value = -537
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H40").Value = 3192.0833
$ws.Range("I40").Value = 3128.4285
$ws.Range("J40").Value = 3281.2
$ws.Range("K40").Value = 3128.4285
$ws.Range("L40").Value = 3281.2
$ws.Range("M40").Value = -2992.4285
$ws.Range("N40").Value = -3553.2
$ws.Range("H46").Value = 2682.15
$ws.Range("I46").Value = 1022.3333
$ws.Range("J46").Value = 4040.182
$ws.Range("K46").Value = 1022.3333
$ws.Range("L46").Value = 4040.182
$ws.Range("M46").Value = -834.3333
$ws.Range("N46").Value = -4416.182
$ws.Range("H122").Value = 3906
$ws.Range("I122").Value = 2303.375
$ws.Range("K122").Value = 6910.125
$ws.Range("M122").Value = -4460.125
$ws.Range("H126").Value = 5828.4287
$ws.Range("I126").Value = 3701.3333
$ws.Range("K126").Value = 11103.9999
$ws.Range("M126").Value = -8633.999899999999
$ws.Range("H130").Value = 109663.336
$ws.Range("J130").Value = 109663.336
$ws.Range("L130").Value = 109663.336
$ws.Range("N130").Value = -119703.336
$ws.Range("H141").Value = 122262.29
$ws.Range("J141").Value = 122262.29
$ws.Range("L141").Value = 122262.29
$ws.Range("N141").Value = -132622.29

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8099.6
$ws.Range("I122").Value = 5364.838
$ws.Range("J122").Value = 13721.056
$ws.Range("K122").Value = 16094.514
$ws.Range("L122").Value = 41163.16800000001
$ws.Range("M122").Value = -13644.514
$ws.Range("N122").Value = -46063.16800000001
$ws.Range("H132").Value = 5078
$ws.Range("I132").Value = 6045.647
$ws.Range("J132").Value = 2728
$ws.Range("K132").Value = 18136.941
$ws.Range("L132").Value = 8184
$ws.Range("M132").Value = -15606.941
$ws.Range("N132").Value = -13244
$ws.Range("H137").Value = 137499.5
$ws.Range("J137").Value = 137499.5
$ws.Range("L137").Value = 137499.5
$ws.Range("N137").Value = -147699.5
